$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Fill in row 38 with the new activity log entry
$ws.Cells.Item(38, 2).Value = 6977
$ws.Cells.Item(38, 3).Value = 43925
$ws.Cells.Item(38, 4).Value = 0.97916666666666663
$ws.Cells.Item(38, 5).Value = 0.98749999999999993
$ws.Cells.Item(38, 7).Value = "Updated folder structure of Documentation directory. Started working on the report."

# Update the selected cell on the sheet to A38
$ws.Range("A38").Select()
